$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 799.8
$ws.Range("I28").Value = 527.4545000000001
$ws.Range("K28").Value = 527.4545000000001
$ws.Range("M28").Value = -42.45450000000005
$ws.Range("H123").Value = 84996.664
$ws.Range("J123").Value = 84996.664
$ws.Range("L123").Value = 84996.664
$ws.Range("N123").Value = -94796.664
$ws.Range("H124").Value = 124000
$ws.Range("J124").Value = 124000
$ws.Range("L124").Value = 124000
$ws.Range("N124").Value = -133820
$ws.Range("H126").Value = 74498.336
$ws.Range("I126").Value = 74750
$ws.Range("J126").Value = 73995
$ws.Range("K126").Value = 74750
$ws.Range("L126").Value = 73995
$ws.Range("M126").Value = -69810
$ws.Range("N126").Value = -83875
$ws.Range("H127").Value = 597.4
$ws.Range("I127").Value = 597.4
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1792.2
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = 3167.8
$ws.Range("H129").Value = 1176.1666
$ws.Range("I129").Value = 1176.1666
$ws.Range("K129").Value = 3528.4998
$ws.Range("M129").Value = 1471.5002
$ws.Range("H138").Value = 2650.055
$ws.Range("J138").Value = 3103.9014
$ws.Range("L138").Value = 9311.7042
$ws.Range("N138").Value = -19591.7042
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 663.5
$ws.Range("I3").Value = 727.25
$ws.Range("J3").Value = 599.75
$ws.Range("K3").Value = 727.25
$ws.Range("L3").Value = 599.75
$ws.Range("M3").Value = -612.25
$ws.Range("N3").Value = -829.75
$ws.Range("H32").Value = 8482838
$ws.Range("I32").Value = 13516633
$ws.Range("K32").Value = 13516633
$ws.Range("M32").Value = -13516346
$ws.Range("H61").Value = 35719644
$ws.Range("I61").Value = 29415196
$ws.Range("K61").Value = 29415196
$ws.Range("M61").Value = -29414984
$ws.Range("H63").Value = 5758.2085
$ws.Range("J63").Value = 7969.231
$ws.Range("L63").Value = 7969.231
$ws.Range("N63").Value = -9341.231
$ws.Range("H66").Value = 5758.2085
$ws.Range("J66").Value = 7969.231
$ws.Range("L66").Value = 39846.155
$ws.Range("N66").Value = -46710.155
$ws.Range("H102").Value = 5512.2666
$ws.Range("I102").Value = 5512.2666
$ws.Range("K102").Value = 5512.2666
$ws.Range("M102").Value = -3890.2666
$ws.Range("H114").Value = 63333
$ws.Range("J114").Value = 63333
$ws.Range("L114").Value = 63333
$ws.Range("N114").Value = -72011
$ws.Range("H132").Value = 10106178
$ws.Range("I132").Value = 15154023
$ws.Range("K132").Value = 45462069
$ws.Range("M132").Value = -45459539
$ws.Range("H133").Value = 75159.336
$ws.Range("J133").Value = 75159.336
$ws.Range("L133").Value = 75159.336
$ws.Range("N133").Value = -80219.336
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 35719644
$ws.Range("I136").Value = 29415196
$ws.Range("K136").Value = 88245588
$ws.Range("M136").Value = -88243038
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 49160
$ws.Range("J116").Value = 49160
$ws.Range("L116").Value = 49160
$ws.Range("N116").Value = -58338
$ws.Range("H31").Value = 1474095.2
$ws.Range("I31").Value = 32719
$ws.Range("J31").Value = 2129266.2
$ws.Range("K31").Value = 32719
$ws.Range("L31").Value = 2129266.2
$ws.Range("M31").Value = -32424
$ws.Range("N31").Value = -2129856.2
$ws.Range("H34").Value = 1474095.2
$ws.Range("I34").Value = 32719
$ws.Range("J34").Value = 2129266.2
$ws.Range("K34").Value = 32719
$ws.Range("L34").Value = 2129266.2
$ws.Range("M34").Value = -32517
$ws.Range("N34").Value = -2129670.2
$ws.Range("H95").Value = 22144.5
$ws.Range("J95").Value = 22144.5
$ws.Range("L95").Value = 22144.5
$ws.Range("N95").Value = -27636.5
$ws.Range("H99").Value = 2719.75
$ws.Range("J99").Value = 3676
$ws.Range("L99").Value = 3676
$ws.Range("N99").Value = -6672
$ws.Range("H126").Value = 2719.75
$ws.Range("J126").Value = 3676
$ws.Range("L126").Value = 11028
$ws.Range("N126").Value = -15968
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2417.8857
$ws.Range("I134").Value = 1310.0435
$ws.Range("K134").Value = 3930.1305
$ws.Range("M134").Value = -1395.1305
$ws.Range("H15").Value = 192.97368
$ws.Range("I15").Value = 163.5
$ws.Range("J15").Value = 203.5
$ws.Range("K15").Value = 490.5
$ws.Range("L15").Value = 610.5
$ws.Range("M15").Value = -350.5
$ws.Range("N15").Value = -890.5
$ws.Range("H60").Value = 1804.1666
$ws.Range("I60").Value = 392.85715
$ws.Range("J60").Value = 3780
$ws.Range("K60").Value = 1178.57145
$ws.Range("L60").Value = 11340
$ws.Range("M60").Value = -927.5714499999999
$ws.Range("N60").Value = -11842
$ws.Range("H113").Value = 1474.1111
$ws.Range("I113").Value = 859.5
$ws.Range("J113").Value = 1965.8
$ws.Range("K113").Value = 2578.5
$ws.Range("L113").Value = 5897.4
$ws.Range("M113").Value = -408.5
$ws.Range("N113").Value = -10237.4
$ws.Range("H131").Value = 3238.5312
$ws.Range("J131").Value = 1709.75
$ws.Range("L131").Value = 5129.25
$ws.Range("N131").Value = -15209.25
$ws.Range("H134").Value = 5001.7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 76.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H99").Value = 72332.664
$ws.Range("I99").Value = 1998
$ws.Range("J99").Value = 107500
$ws.Range("K99").Value = 1998
$ws.Range("L99").Value = 107500
$ws.Range("M99").Value = 248
$ws.Range("N99").Value = -111992
$ws.Range("H107").Value = 533.5
$ws.Range("I107").Value = 253.6
$ws.Range("K107").Value = 253.6
$ws.Range("M107").Value = 1666.4
$ws.Range("H113").Value = 3640.4
$ws.Range("I113").Value = 2636.875
$ws.Range("K113").Value = 2636.875
$ws.Range("M113").Value = -466.875
$ws.Range("H132").Value = 15627234
$ws.Range("I132").Value = 18520702
$ws.Range("K132").Value = 55562106
$ws.Range("M132").Value = -55559576
$ws.Range("H136").Value = 7020.5454
$ws.Range("J136").Value = 7020.5454
$ws.Range("L136").Value = 21061.6362
$ws.Range("N136").Value = -26161.6362
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 558313.9
$ws.Range("I7").Value = 8019.857
$ws.Range("K7").Value = 8019.857
$ws.Range("M7").Value = -7907.857
$ws.Range("H61").Value = 1511.2667
$ws.Range("I61").Value = 1077.1
$ws.Range("K61").Value = 1077.1
$ws.Range("M61").Value = -875.0999999999999
$ws.Range("H81").Value = 109974.5
$ws.Range("J81").Value = 109974.5
$ws.Range("L81").Value = 109974.5
$ws.Range("N81").Value = -111970.5
$ws.Range("H84").Value = 109974.5
$ws.Range("J84").Value = 109974.5
$ws.Range("L84").Value = 329923.5
$ws.Range("N84").Value = -339907.5
$ws.Range("H113").Value = 1511.2667
$ws.Range("I113").Value = 1077.1
$ws.Range("K113").Value = 1077.1
$ws.Range("M113").Value = 1092.9
$ws.Range("H126").Value = 558313.9
$ws.Range("I126").Value = 8019.857
$ws.Range("K126").Value = 24059.571
$ws.Range("M126").Value = -21589.571
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 149.5
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 295
$ws.Range("K17").Value = 4
$ws.Range("L17").Value = 295
$ws.Range("M17").Value = 168
$ws.Range("N17").Value = -639
$ws.Range("H80").Value = 50997.25
$ws.Range("J80").Value = 72000
$ws.Range("L80").Value = 72000
$ws.Range("N80").Value = -73996
$ws.Range("H81").Value = 7471.7144
$ws.Range("I81").Value = 4043.8125
$ws.Range("K81").Value = 8087.625
$ws.Range("M81").Value = -7026.625
$ws.Range("H83").Value = 50997.25
$ws.Range("J83").Value = 72000
$ws.Range("L83").Value = 216000
$ws.Range("N83").Value = -225984
$ws.Range("H84").Value = 7471.7144
$ws.Range("I84").Value = 4043.8125
$ws.Range("K84").Value = 40438.125
$ws.Range("M84").Value = -35134.125
$ws.Range("H113").Value = 590.1177
$ws.Range("J113").Value = 881.875
$ws.Range("L113").Value = 2645.625
$ws.Range("N113").Value = -6985.625

Write-Output "applied sets=213 deletes=2"
